$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.577.09"
$ws.Range("E2").Value = "  +1.22%  "

$ws.Range("D3").Value = "2.489.85"
$ws.Range("E3").Value = "  +0.27%  "

$ws.Range("E4").Value = "  +0.10%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "491.43"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.06%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "152.12"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +7.90%  "

$ws.Range("E7").Value = "  +0.08%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.513"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.31%  "

$ws.Range("D9").Value = "2.499.10"
$ws.Range("E9").Value = "  -0.47%  "

$ws.Range("E10").Value = "  +3.15%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0986"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.46%  "

$ws.Range("E12").Value = "  +0.72%  "

$ws.Range("E13").Value = "  +0.68%  "

$ws.Range("D14").Value = "2.921.83"
$ws.Range("E14").Value = "  -0.24%  "

$ws.Range("D15").Value = "56.720.91"
$ws.Range("E15").Value = "  +1.64%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "21.23"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.73%  "

$ws.Range("E17").Value = "  -1.72%  "

$ws.Range("D18").Value = "2.497.58"
$ws.Range("E18").Value = "  -0.27%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "4.56"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +3.23%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "10.33"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.68%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "321.67"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.49%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.12%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.89"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.94%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "58.42"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.14%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.411"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.52%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.19%  "

$ws.Range("E27").Value = "  -5.92%  "

$ws.Range("D28").Value = "2.593.07"
$ws.Range("E28").Value = "  -0.96%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "7.58"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.55%  "

$ws.Range("D30").Value = "0.0₃0803"
$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("E31").Value = "  +0.13%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "151.59"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.69%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "18.31"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.10%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.52"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.94%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.26"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.35%  "

$ws.Range("E36").Value = "  +2.69%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.78"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.59%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.871"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.13%  "

$ws.Range("E39").Value = "  +4.77%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "34.12"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.88%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "3.53"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +2.50%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.0563"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.01%  "

$ws.Range("E43").Value = "  -0.23%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.995"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "266.16"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.28%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "4.81"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.06%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0931"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +1.92%  "

$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "10.23"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.99%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0229"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.38%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "17.80"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +1.26%  "

$ws.Range("D51").Value = "1.892.40"
$ws.Range("E51").Value = "  -5.61%  "
